$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.743.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.182.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '402.36'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.56'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.67%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.678.36'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.13'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.56%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.06'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.07'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +7.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.188.29'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.60'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '54.605.74'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.33'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.97'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.92%  '
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000100'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.32'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '276.66'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.19'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.86'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.50'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.29'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0507'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +8.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.04'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.62%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.59'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.94'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +13.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.07'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '134.13'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.49%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.93'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.291'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.03'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.118'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.31'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.46'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +39.04%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.102.69'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0511'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.30%  '
